$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) column C for rows 2..222 from 45192 to 45202
for ($r = 2; $r -le 222; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# Row 222 needs explicit row height (matches formatting of surrounding rows)
$ws.Rows.Item(222).RowHeight = 15

# Add new row 223 with data
$ws.Cells.Item(223, 1).Value = "A 47057-2023"
$ws.Cells.Item(223, 2).Value = 45196
$ws.Cells.Item(223, 3).Value = 45202
$ws.Cells.Item(223, 4).Value = "VÄSTRA GÖTALANDS LÄN"
$ws.Cells.Item(223, 5).Value = "HERRLJUNGA"
$ws.Cells.Item(223, 7).Value = 0.7
$ws.Cells.Item(223, 8).Value = 0
$ws.Cells.Item(223, 9).Value = 0
$ws.Cells.Item(223, 10).Value = 0
$ws.Cells.Item(223, 11).Value = 0
$ws.Cells.Item(223, 12).Value = 0
$ws.Cells.Item(223, 13).Value = 0
$ws.Cells.Item(223, 14).Value = 0
$ws.Cells.Item(223, 15).Value = 0
$ws.Cells.Item(223, 16).Value = 0
$ws.Cells.Item(223, 17).Value = 0

# Match formatting of the other rows for the new row
$ws.Cells.Item(223, 2).NumberFormat = $ws.Cells.Item(222, 2).NumberFormat
$ws.Cells.Item(223, 3).NumberFormat = $ws.Cells.Item(222, 3).NumberFormat
$ws.Cells.Item(223, 18).WrapText = $true

